$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "66.342.84"
Set-TextValue "E2" "  -0.38%  "
Set-TextValue "D3" "3.313.29"
Set-TextValue "E3" "  -1.82%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.13%  "
Set-TextValue "D5" "190.06"
Set-TextValue "E5" "  +3.59%  "
Set-TextValue "D6" "562.48"
Set-TextValue "E6" "  +0.25%  "
Set-TextValue "E7" "  -0.12%  "
Set-TextValue "D8" "0.590"
Set-TextValue "E8" "  -1.54%  "
Set-TextValue "D9" "3.305.41"
Set-TextValue "E9" "  -1.77%  "
Set-TextValue "E10" "  -1.15%  "
Set-TextValue "D11" "0.589"
Set-TextValue "E11" "  -1.22%  "
Set-TextValue "D12" "48.04"
Set-TextValue "E12" "  -0.26%  "
Set-TextValue "E13" "  +1.23%  "
Set-TextValue "E14" "  -0.32%  "
Set-TextValue "D15" "3.841.95"
Set-TextValue "E15" "  -1.90%  "
Set-TextValue "D16" "614.56"
Set-TextValue "E16" "  +1.58%  "
Set-TextValue "D17" "18.13"
Set-TextValue "E17" "  -0.71%  "
Set-TextValue "D18" "66.327.30"
Set-TextValue "E18" "  -0.26%  "
Set-TextValue "E19" "  -0.17%  "
Set-TextValue "D20" "3.302.68"
Set-TextValue "E20" "  -2.29%  "
Set-TextValue "E21" "  -4.09%  "
Set-TextValue "D22" "0.914"
Set-TextValue "E22" "  -0.43%  "
Set-TextValue "D23" "18.66"
Set-TextValue "E23" "  +9.75%  "
Set-TextValue "D24" "5.13"
Set-TextValue "E24" "  -1.29%  "
Set-TextValue "D25" "101.83"
Set-TextValue "E25" "  +3.57%  "
Set-TextValue "E26" "  -1.77%  "
Set-TextValue "E27" "  +0.02%  "
Set-TextValue "D28" "2.77"
Set-TextValue "E28" "  +1.40%  "
Set-TextValue "D29" "9.87"
Set-TextValue "E29" "  +3.45%  "
Set-TextValue "D30" "8.66"
Set-TextValue "E30" "  -1.68%  "
Set-TextValue "D31" "30.46"
Set-TextValue "E31" "  -1.09%  "
Set-TextValue "D32" "6.80"
Set-TextValue "E32" "  +7.40%  "
Set-TextValue "D33" "4.10"
Set-TextValue "E33" "  +7.01%  "
Set-TextValue "D34" "576.83"
Set-TextValue "E34" "  +3.74%  "
Set-TextValue "D35" "11.15"
Set-TextValue "E35" "  -0.33%  "
Set-TextValue "E36" "  -0.21%  "
Set-TextValue "D37" "3.735.38"
Set-TextValue "E37" "  -3.23%  "
Set-TextValue "B38" "Dai"
Set-TextValue "C38" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D38" "1.00"
Set-TextValue "E38" "  +0.31%  "
Set-TextValue "B39" "OKB"
Set-TextValue "C39" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D39" "57.29"
Set-TextValue "E39" "  -1.43%  "
Set-TextValue "D40" "0.0₃0734"
Set-TextValue "E40" "  +1.27%  "
Set-TextValue "E41" "  -1.90%  "
Set-TextValue "D42" "3.49"
Set-TextValue "E42" "  +1.13%  "
Set-TextValue "D43" "34.18"
Set-TextValue "E43" "  +4.85%  "
Set-TextValue "D44" "0.131"
Set-TextValue "E44" "  +1.59%  "
Set-TextValue "E45" "  +1.32%  "
Set-TextValue "E46" "  -2.48%  "
Set-TextValue "E47" "  +1.74%  "
Set-TextValue "E48" "  +2.16%  "
Set-TextValue "D49" "0.130"
Set-TextValue "E49" "  -0.93%  "
Set-TextValue "D50" "2.60"
Set-TextValue "E50" "  -3.16%  "
Set-TextValue "D51" "0.999"
Set-TextValue "E51" "  -0.04%  "
